# chore: update Sheets via scheduled runner
# Refresh market-price-derived columns (currentAveragePrice / NQ / HQ,
# LevePriceNQ / HQ, LeveProfitNQ / HQ) for a handful of leve rows across
# the eight crafting-job sheets, reflecting newer Universalis price data.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 507.80435
$ws.Range("J17").Value = 298.92105
$ws.Range("L17").Value = 896.76315
$ws.Range("N17").Value = -1232.76315
$ws.Range("H112").Value = 1267.9286
$ws.Range("J112").Value = 1325.0769
$ws.Range("L112").Value = 3975.2307
$ws.Range("N112").Value = -6191.2307
$ws.Range("H129").Value = 812.76
$ws.Range("I129").Value = 249
$ws.Range("J129").Value = 861.7826
$ws.Range("K129").Value = 747
$ws.Range("L129").Value = 2585.3478
$ws.Range("M129").Value = 4253
$ws.Range("N129").Value = -12585.3478
$ws.Range("H137").Value = 1288961.1
$ws.Range("I137").Value = 1588652.4
$ws.Range("J137").Value = 4570
$ws.Range("K137").Value = 4765957.199999999
$ws.Range("L137").Value = 13710
$ws.Range("M137").Value = -4763407.199999999
$ws.Range("N137").Value = -18810
$ws.Range("H138").Value = 2481.77
$ws.Range("I138").Value = 625.85
$ws.Range("J138").Value = 2945.75
$ws.Range("K138").Value = 1877.55
$ws.Range("L138").Value = 8837.25
$ws.Range("M138").Value = 3262.45
$ws.Range("N138").Value = -19117.25
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 6109.375
$ws.Range("I32").Value = 4586.2617
$ws.Range("K32").Value = 4586.2617
$ws.Range("M32").Value = -4299.2617
$ws.Range("H53").Value = 0
$ws.Range("J53").Value = 0
$ws.Range("L53").Value = 0
$ws.Range("N53").Value = ""
$ws.Range("H74").Value = 1834.84
$ws.Range("I74").Value = 1146.6666
$ws.Range("K74").Value = 1146.6666
$ws.Range("M74").Value = -272.6666
$ws.Range("H77").Value = 1834.84
$ws.Range("I77").Value = 1146.6666
$ws.Range("K77").Value = 5733.333000000001
$ws.Range("M77").Value = -1365.333000000001
$ws.Range("H122").Value = 8000
$ws.Range("I122").Value = 1000
$ws.Range("K122").Value = 3000
$ws.Range("M122").Value = -550
$ws.Range("H137").Value = 43780
$ws.Range("J137").Value = 43780
$ws.Range("L137").Value = 43780
$ws.Range("N137").Value = -53980
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H35").Value = 25251.2
$ws.Range("J35").Value = 25251.2
$ws.Range("L35").Value = 25251.2
$ws.Range("N35").Value = -25871.2
$ws.Range("H82").Value = 22333.666
$ws.Range("I82").Value = 4039.25
$ws.Range("J82").Value = 36969.2
$ws.Range("K82").Value = 4039.25
$ws.Range("L82").Value = 36969.2
$ws.Range("M82").Value = -3656.25
$ws.Range("N82").Value = -37735.2
$ws.Range("H85").Value = 22333.666
$ws.Range("I85").Value = 4039.25
$ws.Range("J85").Value = 36969.2
$ws.Range("K85").Value = 4039.25
$ws.Range("L85").Value = 36969.2
$ws.Range("M85").Value = -2713.25
$ws.Range("N85").Value = -39621.2
$ws.Range("H94").Value = 1063.8948
$ws.Range("I94").Value = 1033.6666
$ws.Range("J94").Value = 1177.25
$ws.Range("K94").Value = 1033.6666
$ws.Range("L94").Value = 1177.25
$ws.Range("M94").Value = -582.6666
$ws.Range("N94").Value = -2079.25
$ws.Range("H137").Value = 40032.5
$ws.Range("J137").Value = 40032.5
$ws.Range("L137").Value = 40032.5
$ws.Range("N137").Value = -50232.5
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 1292.8572
$ws.Range("I16").Value = 1170.3334
$ws.Range("J16").Value = 1384.75
$ws.Range("K16").Value = 1170.3334
$ws.Range("L16").Value = 1384.75
$ws.Range("M16").Value = -883.3334
$ws.Range("N16").Value = -1958.75
$ws.Range("H31").Value = 8993.200000000001
$ws.Range("I31").Value = 1548.2
$ws.Range("K31").Value = 1548.2
$ws.Range("M31").Value = -1253.2
$ws.Range("H34").Value = 8993.200000000001
$ws.Range("I34").Value = 1548.2
$ws.Range("K34").Value = 1548.2
$ws.Range("M34").Value = -1346.2
$ws.Range("H99").Value = 9094855
$ws.Range("I99").Value = 15386831
$ws.Range("J99").Value = 6444.4443
$ws.Range("K99").Value = 15386831
$ws.Range("L99").Value = 6444.4443
$ws.Range("M99").Value = -15385333
$ws.Range("N99").Value = -9440.444299999999
$ws.Range("H113").Value = 1292.8572
$ws.Range("I113").Value = 1170.3334
$ws.Range("J113").Value = 1384.75
$ws.Range("K113").Value = 1170.3334
$ws.Range("L113").Value = 1384.75
$ws.Range("M113").Value = 999.6666
$ws.Range("N113").Value = -5724.75
$ws.Range("H126").Value = 9094855
$ws.Range("I126").Value = 15386831
$ws.Range("J126").Value = 6444.4443
$ws.Range("K126").Value = 46160493
$ws.Range("L126").Value = 19333.3329
$ws.Range("M126").Value = -46158023
$ws.Range("N126").Value = -24273.3329
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H58").Value = 2488.889
$ws.Range("I58").Value = 1566.6666
$ws.Range("J58").Value = 4333.3335
$ws.Range("K58").Value = 4699.9998
$ws.Range("L58").Value = 13000.0005
$ws.Range("M58").Value = -4571.9998
$ws.Range("N58").Value = -13256.0005
$ws.Range("H95").Value = 5000
$ws.Range("J95").Value = 5000
$ws.Range("L95").Value = 15000
$ws.Range("N95").Value = -19118
$ws.Range("H113").Value = 2976867.5
$ws.Range("I113").Value = 578.9583
$ws.Range("J113").Value = 6945252.5
$ws.Range("K113").Value = 1736.8749
$ws.Range("L113").Value = 20835757.5
$ws.Range("M113").Value = 433.1251
$ws.Range("N113").Value = -20840097.5
$ws.Range("H131").Value = 852.6083
$ws.Range("J131").Value = 874.27954
$ws.Range("L131").Value = 2622.83862
$ws.Range("N131").Value = -12702.83862
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H46").Value = 16643.143
$ws.Range("J46").Value = 16643.143
$ws.Range("L46").Value = 16643.143
$ws.Range("N46").Value = -16955.143
$ws.Range("H137").Value = 40273.2
$ws.Range("J137").Value = 40273.2
$ws.Range("L137").Value = 40273.2
$ws.Range("N137").Value = -50473.2
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H53").Value = 13999.5
$ws.Range("J53").Value = 13999.5
$ws.Range("L53").Value = 13999.5
$ws.Range("N53").Value = -15035.5
$ws.Range("H96").Value = 38800
$ws.Range("J96").Value = 38800
$ws.Range("L96").Value = 38800
$ws.Range("N96").Value = -44292
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 13895301
$ws.Range("I132").Value = 8267.357
$ws.Range("J132").Value = 33337148
$ws.Range("K132").Value = 24802.071
$ws.Range("L132").Value = 100011444
$ws.Range("M132").Value = -22272.071
$ws.Range("N132").Value = -100016504
$ws.Range("H136").Value = 6164.162
$ws.Range("I136").Value = 5488.36
$ws.Range("J136").Value = 7572.0835
$ws.Range("K136").Value = 16465.08
$ws.Range("L136").Value = 22716.2505
$ws.Range("M136").Value = -13915.08
$ws.Range("N136").Value = -27816.2505
